$wb = $excel.ActiveWorkbook
Write-Host "noop"
